# Updates cryptocurrency price/volume data on the "cryptos" worksheet.
# Generated to reflect the scraped price/percentage refresh described in the
# commit "Updated cryptos list ... with GitHub Actions".

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of row number -> (new Price value, new Volume(1h) value, forceTextForPrice)
$data = @{
    2 = @("69.346.54", "  -1.88%  ", 0)
    3 = @("3.492.43", "  -1.88%  ", 0)
    4 = @("0.998", "  -0.16%  ", 1)
    5 = @("611.04", "  +4.89%  ", 1)
    6 = @("185.80", "  +0.29%  ", 1)
    7 = @("0.634", "  +0.03%  ", 1)
    8 = @($null, "  -0.08%  ", 0)
    9 = @($null, "  -0.10%  ", 0)
    10 = @("0.653", "  -0.02%  ", 1)
    11 = @("52.98", "  -2.67%  ", 1)
    12 = @("0.0000312", "  -1.04%  ", 1)
    13 = @("9.58", "  +1.00%  ", 1)
    14 = @("4.048.25", "  -1.82%  ", 0)
    15 = @("600.64", "  +5.50%  ", 1)
    16 = @("69.350.82", "  -1.88%  ", 0)
    17 = @("18.88", "  -1.92%  ", 1)
    18 = @("12.59", "  +1.93%  ", 1)
    19 = @("3.504.37", "  -1.47%  ", 0)
    20 = @($null, "  -0.21%  ", 0)
    21 = @($null, "  -1.40%  ", 0)
    22 = @("17.23", "  -3.01%  ", 1)
    23 = @("105.21", "  +10.29%  ", 1)
    24 = @("4.67", "  +2.76%  ", 1)
    25 = @("5.06", "  +1.65%  ", 1)
    26 = @($null, "  +3.01%  ", 0)
    27 = @("10.97", "  -3.11%  ", 1)
    28 = @("9.98", "  +9.30%  ", 1)
    29 = @("33.57", "  +3.35%  ", 1)
    30 = @("6.98", "  -3.49%  ", 1)
    31 = @("12.41", "  +1.10%  ", 1)
    32 = @($null, "  -0.18%  ", 0)
    33 = @("3.90", "  +16.61%  ", 1)
    34 = @("63.29", "  -0.34%  ", 1)
    35 = @($null, "  -7.90%  ", 0)
    36 = @($null, "  -0.08%  ", 0)
    37 = @("520.09", "  -4.71%  ", 1)
    38 = @($null, "  -4.07%  ", 0)
    39 = @("3.610.18", "  +1.02%  ", 0)
    40 = @("3.61", "  +5.22%  ", 1)
    41 = @("36.81", "  -2.51%  ", 1)
    42 = @("0.0₃0780", "  -1.63%  ", 0)
    43 = @("0.139", "  +1.55%  ", 1)
    44 = @($null, "  +2.78%  ", 0)
    45 = @($null, "  +1.27%  ", 0)
    46 = @($null, "  +3.46%  ", 0)
    47 = @("3.33", "  -4.36%  ", 1)
    48 = @("8.81", "  -5.71%  ", 1)
    49 = @($null, "  +0.34%  ", 0)
    50 = @("0.000246", "  -6.00%  ", 1)
    51 = @($null, "  -9.93%  ", 0)
}

foreach ($row in $data.Keys) {
    $entry = $data[$row]
    $dVal = $entry[0]
    $eVal = $entry[1]
    $forceText = $entry[2]

    if ($dVal -ne $null) {
        $dCell = $ws.Cells.Item($row, 4)
        if ($forceText -eq 1) {
            # These values parse as plain numbers (e.g. "0.998"); force them to
            # be stored as text so the displayed price string is preserved
            # exactly (matching the original inline-string formatting), then
            # restore the default "Normal" style so no stray number format is
            # left behind on the cell.
            $dCell.NumberFormat = "@"
            $dCell.Value = $dVal
            $dCell.Style = "Normal"
        } else {
            $dCell.Value = $dVal
        }
    }

    if ($eVal -ne $null) {
        $ws.Cells.Item($row, 5).Value = $eVal
    }
}
